# Daily update at 8 AM UTC
# Appends the next day's row (row 99) to the "Wins Over Time" tracking sheet
# and moves the "last row" date formatting from the old last row (98) to the
# new last row (99).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 98 was previously the last data row and used the special "last row"
# date format (YYYY-MM-DD, no time). Since it is no longer the last row,
# give it the regular date/time number format used by all the other
# non-final data rows.
$ws.Range("A98").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Add the new day's data in row 99.
$ws.Range("A99").Value = 45686
$ws.Range("B99").Value = 239
$ws.Range("C99").Value = 229
$ws.Range("D99").Value = 231

# The new row is now the last row, so it gets the special "last row" date
# format that row 98 used to have.
$ws.Range("A99").NumberFormat = "YYYY-MM-DD"
